# Staff dashboard / book list touch-ups:
#  - Update check-out status, who-checked, and expected-return date for a
#    handful of existing rows.
#  - Append a batch of newly catalogued books as rows 16-40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style for "Expected Return" date cells (style index 2 / custom
# date number format) - copy it from an existing formatted cell (H2) so new
# / retouched date cells line up with the rest of the column.
$dateFormat = $ws.Cells.Item(2, 8).NumberFormat

# ---------------------------------------------------------------------
# 1. Row 7 - The Great Gatsby: now checked out by Galactus, new due date.
# ---------------------------------------------------------------------
$ws.Cells.Item(7, 6).Value = "no"
$ws.Cells.Item(7, 7).Value = "Galactus"
$ws.Cells.Item(7, 8).Value = 45883
$ws.Cells.Item(7, 8).NumberFormat = $dateFormat

# ---------------------------------------------------------------------
# 2. Row 8 - One Hundred Years of Solitude: returned (back in stock).
# ---------------------------------------------------------------------
$ws.Cells.Item(8, 6).Value = "Yes"
$ws.Cells.Item(8, 7).Value = ""

# ---------------------------------------------------------------------
# 3. Row 9 - Wildblood: normalize casing, corrected due date.
# ---------------------------------------------------------------------
$ws.Cells.Item(9, 6).Value = "no"
$ws.Cells.Item(9, 8).Value = 45872

# ---------------------------------------------------------------------
# 4. Row 10 - Algebra of Happiness: now checked out by BillGates.
# ---------------------------------------------------------------------
$ws.Cells.Item(10, 6).Value = "No"
$ws.Cells.Item(10, 7).Value = "BillGates"
$ws.Cells.Item(10, 8).Value = 45881

# ---------------------------------------------------------------------
# 5. Row 13 - Who Moved My Cheese?: now checked out by admin.
# ---------------------------------------------------------------------
$ws.Cells.Item(13, 6).Value = "no"
$ws.Cells.Item(13, 7).Value = "admin"
$ws.Cells.Item(13, 8).Value = 45881
$ws.Cells.Item(13, 8).NumberFormat = $dateFormat

# ---------------------------------------------------------------------
# 6. Append newly catalogued books as rows 16-40.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=16; Title='The Night Circus'; Author='Erin Morgenstern'; Genre='Fantasy'; ISBN=9780385534635; Year=2011; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=17; Title='Sapiens'; Author='Yuval Noah Harari'; Genre='Non-Fiction'; ISBN=9780062316097; Year=2015; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=18; Title='Circe'; Author='Madeline Miller'; Genre='Mythology'; ISBN=9780316556323; Year=2018; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=19; Title='Educated'; Author='Tara Westover'; Genre='Memoir'; ISBN=9780399590504; Year=2018; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=20; Title='The Midnight Library'; Author='Matt Haig'; Genre='Fantasy'; ISBN=9780525559474; Year=2020; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=21; Title='Pachinko'; Author='Min Jin Lee'; Genre='Historical'; ISBN=9781455563920; Year=2017; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=22; Title='The Martian'; Author='Andy Weir'; Genre='Science Fiction'; ISBN=9780553418026; Year=2014; InStock='no'; Who='Galactus'; DueDate=45881 }
    @{ Row=23; Title='The Alchemist'; Author='Paulo Coelho'; Genre='Philosophical'; ISBN=9780061122415; Year=1993; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=24; Title='Normal People'; Author='Sally Rooney'; Genre='Drama'; ISBN=9781984822178; Year=2019; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=25; Title='Project Hail Mary'; Author='Andy Weir'; Genre='Science Fiction'; ISBN=9780593135204; Year=2021; InStock='no'; Who='Bowser'; DueDate=45881 }
    @{ Row=26; Title='Where the Crawdads Sing'; Author='Delia Owens'; Genre='Mystery'; ISBN=9780735219106; Year=2018; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=27; Title='A Man Called Ove'; Author='Fredrik Backman'; Genre='Drama'; ISBN=9781476738024; Year=2012; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=28; Title='The Silent Patient'; Author='Alex Michaelides'; Genre='Thriller'; ISBN=9781250301697; Year=2019; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=29; Title='The Song of Achilles'; Author='Madeline Miller'; Genre='Mythology'; ISBN=9780062060624; Year=2012; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=30; Title='Daisy Jones & The Six'; Author='Taylor Jenkins Reid'; Genre='Fiction'; ISBN=9781524798628; Year=2019; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=31; Title='Klara and the Sun'; Author='Kazuo Ishiguro'; Genre='Science Fiction'; ISBN=9780593318171; Year=2021; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=32; Title='Mexican Gothic'; Author='Silvia Moreno-Garcia'; Genre='Horror'; ISBN=9780525620785; Year=2020; InStock='no'; Who='JohnDoe'; DueDate=45856 }
    @{ Row=33; Title='The Paris Library'; Author='Janet Skeslien Charles'; Genre='Historical'; ISBN=9781982134198; Year=2021; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=34; Title='The Guest List'; Author='Lucy Foley'; Genre='Mystery'; ISBN=9780062868930; Year=2020; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=35; Title='The Seven Husbands of Evelyn Hugo'; Author='Taylor Jenkins Reid'; Genre='Drama'; ISBN=9781501161933; Year=2017; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=36; Title='Beach Read'; Author='Emily Henry'; Genre='Romance'; ISBN=9781984806734; Year=2020; InStock='no'; Who='BillGates'; DueDate=45881 }
    @{ Row=37; Title='Verity'; Author='Colleen Hoover'; Genre='Thriller'; ISBN=9781538724736; Year=2018; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=38; Title='The Book Thief'; Author='Markus Zusak'; Genre='Historical'; ISBN=9780375842207; Year=2005; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=39; Title='Little Fires Everywhere'; Author='Celeste Ng'; Genre='Fiction'; ISBN=9780735224315; Year=2017; InStock='yes'; Who=$null; DueDate=$null }
    @{ Row=40; Title='The Paper Palace'; Author='Miranda Cowley Heller'; Genre='Fiction'; ISBN=9780593329825; Year=2021; InStock='yes'; Who=$null; DueDate=$null }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Title
    $ws.Cells.Item($r.Row, 2).Value = $r.Author
    $ws.Cells.Item($r.Row, 3).Value = $r.Genre
    $ws.Cells.Item($r.Row, 4).Value = $r.ISBN
    $ws.Cells.Item($r.Row, 5).Value = $r.Year
    $ws.Cells.Item($r.Row, 6).Value = $r.InStock
    if ($r.Who) {
        $ws.Cells.Item($r.Row, 7).Value = $r.Who
    }
    if ($r.DueDate) {
        $ws.Cells.Item($r.Row, 8).Value = $r.DueDate
        $ws.Cells.Item($r.Row, 8).NumberFormat = $dateFormat
    }
}
